# Apply NATMI edit: expand LR-pair results from a single row to full
# sender x target cluster combinations (FAPs, sCs) x (FAPs, sCs, ECs),
# following Dr Hou's advice.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("FAPs", "Fgf13", "Scn5a", "ECs", 1, 0.3333333333333333, 0.03538133333333333, 0.106144, 0.0113454327482436, 0.01134543274824361, 2, 0.6666666666666666, 0.082317, 0.246951, 0.05525975897489143, 0.05525975897489144, 0.002912485216, 0.026212366944, 0.0006269458791337816, 0.0006269458791337818),
    @("FAPs", "Fgf13", "Scn5a", "FAPs", 1, 0.3333333333333333, 0.03538133333333333, 0.106144, 0.0113454327482436, 0.01134543274824361, 1, 0.3333333333333333, 0.08236, 0.24708, 0.05528862506131247, 0.05528862506131247, 0.002914006613333334, 0.02622605952, 0.0006272733773759765, 0.0006272733773759766),
    @("FAPs", "Fgf13", "Scn5a", "sCs", 1, 0.3333333333333333, 0.03538133333333333, 0.106144, 0.0113454327482436, 0.01134543274824361, 3, 1, 1.324960333333333, 3.974881, 0.8894516159637962, 0.8894516159637962, 0.04687886320711111, 0.421909768864, 0.01009121349173385, 0.01009121349173385),
    @("sCs", "Fgf13", "Scn5a", "ECs", 3, 1, 3.083171666666666, 9.249514999999999, 0.9886545672517564, 0.9886545672517564, 2, 0.6666666666666666, 0.082317, 0.246951, 0.05525975897489143, 0.05525975897489144, 0.253797442085, 2.284176978765, 0.05463281309575765, 0.05463281309575765),
    @("sCs", "Fgf13", "Scn5a", "FAPs", 3, 1, 3.083171666666666, 9.249514999999999, 0.9886545672517564, 0.9886545672517564, 1, 0.3333333333333333, 0.08236, 0.24708, 0.05528862506131247, 0.05528862506131247, 0.2539300184666666, 2.285370166199999, 0.05466135168393649, 0.05466135168393649),
    @("sCs", "Fgf13", "Scn5a", "sCs", 3, 1, 3.083171666666666, 9.249514999999999, 0.9886545672517564, 0.9886545672517564, 3, 1, 1.324960333333333, 3.974881, 0.8894516159637962, 0.8894516159637962, 4.085080159190555, 36.765721432715, 0.8793604024720624, 0.8793604024720624)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}
